$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).Value = 30
$ws.Cells.Item(11, 2).Value = 3
$ws.Cells.Item(11, 3).Value = "2021-04-14 23:16:00.140322"
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0

$ws.Cells.Item(12, 1).Value = 30
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 3).Value = "2021-04-14 23:16:51.742843"
$ws.Cells.Item(12, 4).Value = 2
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 2
$ws.Cells.Item(12, 9).Value = 2

$ws.Cells.Item(13, 1).Value = 30
$ws.Cells.Item(13, 2).Value = 3
$ws.Cells.Item(13, 3).Value = "2021-04-14 23:46:32.245210"
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 2
$ws.Cells.Item(13, 9).Value = 2

$ws.Cells.Item(14, 1).Value = 30
$ws.Cells.Item(14, 2).Value = 3
$ws.Cells.Item(14, 3).Value = "2021-04-14 23:48:25.206486"
$ws.Cells.Item(14, 4).Value = 3
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 2
$ws.Cells.Item(14, 7).Value = 3
$ws.Cells.Item(14, 8).Value = 3
$ws.Cells.Item(14, 9).Value = 3

$ws.Cells.Item(15, 1).Value = 15
$ws.Cells.Item(15, 2).Value = 3
$ws.Cells.Item(15, 3).Value = "2021-04-14 23:50:16.661164"
$ws.Cells.Item(15, 4).Value = 7
$ws.Cells.Item(15, 5).Value = 7
$ws.Cells.Item(15, 6).Value = 7
$ws.Cells.Item(15, 7).Value = 7
$ws.Cells.Item(15, 8).Value = 7
$ws.Cells.Item(15, 9).Value = 7

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 3
$ws.Cells.Item(16, 3).Value = "2021-04-14 23:54:28.575005"
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 1
$ws.Cells.Item(16, 9).Value = 1

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 3
$ws.Cells.Item(17, 3).Value = "2021-04-15 13:12:56.141677"
$ws.Cells.Item(17, 4).Value = 9
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 9
$ws.Cells.Item(17, 7).Value = 3
$ws.Cells.Item(17, 8).Value = 9
$ws.Cells.Item(17, 9).Value = 9

$ws.Cells.Item(18, 1).Value = 15
$ws.Cells.Item(18, 2).Value = 3
$ws.Cells.Item(18, 3).Value = "2021-04-15 13:45:35.189974"
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(18, 9).Value = 1

$ws.Cells.Item(19, 1).Value = 15
$ws.Cells.Item(19, 2).Value = 3
$ws.Cells.Item(19, 3).Value = "2021-04-15 13:46:01.736827"
$ws.Cells.Item(19, 4).Value = 2
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 2
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 2
$ws.Cells.Item(19, 9).Value = 2

$ws.Cells.Item(20, 1).Value = 15
$ws.Cells.Item(20, 2).Value = 3
$ws.Cells.Item(20, 3).Value = "2021-04-15 13:46:47.675392"
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(20, 8).Value = 1
$ws.Cells.Item(20, 9).Value = 1

$ws.Cells.Item(21, 1).Value = 15
$ws.Cells.Item(21, 2).Value = 3
$ws.Cells.Item(21, 3).Value = "2021-04-15 13:48:23.374307"
$ws.Cells.Item(21, 4).Value = 3
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 4
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(21, 8).Value = 4
$ws.Cells.Item(21, 9).Value = 4

$ws.Cells.Item(22, 1).Value = 15
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(22, 3).Value = "2021-04-15 13:51:07.392178"
$ws.Cells.Item(22, 4).Value = 3
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 3
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 3
$ws.Cells.Item(22, 9).Value = 3

$ws.Cells.Item(23, 1).Value = 15
$ws.Cells.Item(23, 2).Value = 3
$ws.Cells.Item(23, 3).Value = "2021-04-15 14:10:51.898023"
$ws.Cells.Item(23, 4).Value = 2
$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 6).Value = 2
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 2
$ws.Cells.Item(23, 9).Value = 2

$ws.Cells.Item(24, 1).Value = 15
$ws.Cells.Item(24, 2).Value = 3
$ws.Cells.Item(24, 3).Value = "2021-04-15 14:12:56.397366"
$ws.Cells.Item(24, 4).Value = 3
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 3
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = 3
$ws.Cells.Item(24, 9).Value = 3
